$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.582.96'
$ws.Range("E2").Value = '  -0.23%  '

$ws.Range("D3").Value = '1.687.20'
$ws.Range("E3").Value = '  -0.03%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.008'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +0.51%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '313.87'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.53%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '1.010'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.85%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.3900'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.99%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.4022'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -0.57%  '

$ws.Range("B9").Value = 'BinanceUSD'
$ws.Range("C9").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '1.009'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.62%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '1.481'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -0.37%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '52.93'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +0.18%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.08684'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -1.29%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '7.582'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +5.25%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '24.46'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +4.37%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '7.924'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.53%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.00001330'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +1.36%  '

$ws.Range("D17").Value = '1.694.92'
$ws.Range("E17").Value = '  -0.07%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '98.22'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -1.42%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.07107'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +1.52%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '19.60'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.92%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '7.252'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +3.69%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '1.011'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.84%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '14.13'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.87%  '

$ws.Range("D24").Value = '24.579.32'
$ws.Range("E24").Value = '  -0.18%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.999'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -8.70%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.350'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.59%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '22.58'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.25%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '161.33'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.86%  '

$ws.Range("E29").Value = '  +10.98%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '5.242'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +1.41%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '136.19'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.76%  '

$ws.Range("D32").Value = '1.877.50'
$ws.Range("E32").Value = '  -0.21%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '7.509'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +4.72%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.08716'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +2.05%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.026'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -2.82%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '1.981'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +4.76%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.02889'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +6.64%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.2707'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -0.84%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '10.66'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -4.01%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.09093'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.83%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '14.03'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -1.80%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.7693'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +1.55%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.452'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.33%  '

$ws.Range("E44").Value = '  +3.68%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.7101'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -0.19%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.553'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.97%  '

$ws.Range("E47").Value = '  -0.26%  '

$ws.Range("E48").Value = '  +0.40%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '1.326'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +0.95%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '137.83'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -0.99%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '90.53'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +1.31%  '
